# Update the LR-pairs data table (rows 2-17) with recomputed NATMI values,
# following Dr Hou's advice: now uses 3 replicate samples instead of 1, and
# adds "sCs" as an additional target cluster, expanding the table to A1:T17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Jag2", "Notch3", "ECs", 3, 1, 20.70050833333333, 62.101525, 0.8277101186170105, 0.8277101186170105, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 358.6347198360278, 3227.712478524249, 0.1499288810582576, 0.1499288810582576)
    ,@("ECs", "Jag2", "Notch3", "FAPs", 3, 1, 20.70050833333333, 62.101525, 0.8277101186170105, 0.8277101186170105, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 136.4262792604972, 1227.836513344475, 0.05703362855057553, 0.05703362855057553)
    ,@("ECs", "Jag2", "Notch3", "M2", 3, 1, 20.70050833333333, 62.101525, 0.8277101186170105, 0.8277101186170105, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 2.301793024125, 20.716137217125, 0.0009622750766923811, 0.0009622750766923811)
    ,@("ECs", "Jag2", "Notch3", "sCs", 3, 1, 20.70050833333333, 62.101525, 0.8277101186170105, 0.8277101186170105, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1482.546511546547, 13342.91860391892, 0.619785333931485, 0.6197853339314849)
    ,@("FAPs", "Jag2", "Notch3", "ECs", 3, 1, 0.7925996666666667, 2.377799, 0.03169210888681734, 0.03169210888681734, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 13.73172845902556, 123.58555613123, 0.005740611739429007, 0.005740611739429007)
    ,@("FAPs", "Jag2", "Notch3", "FAPs", 3, 1, 0.7925996666666667, 2.377799, 0.03169210888681734, 0.03169210888681734, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 5.223611986973445, 47.012507882761, 0.002183754826212882, 0.002183754826212882)
    ,@("FAPs", "Jag2", "Notch3", "M2", 3, 1, 0.7925996666666667, 2.377799, 0.03169210888681734, 0.03169210888681734, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 0.088133119935, 0.7931980794150001, 0.00003684445293548053, 0.00003684445293548053)
    ,@("FAPs", "Jag2", "Notch3", "sCs", 3, 1, 0.7925996666666667, 2.377799, 0.03169210888681734, 0.03169210888681734, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 56.76507320245145, 510.885658822063, 0.02373089786823997, 0.02373089786823997)
    ,@("M2", "Jag2", "Notch3", "ECs", 3, 1, 1.562510666666667, 4.687532, 0.06247701111592723, 0.06247701111592723, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 27.07037750751556, 243.63339756764, 0.01131689483768356, 0.01131689483768356)
    ,@("M2", "Jag2", "Notch3", "FAPs", 3, 1, 1.562510666666667, 4.687532, 0.06247701111592723, 0.06247701111592723, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 10.29769477761644, 92.67925299854801, 0.004304998289606196, 0.004304998289606196)
    ,@("M2", "Jag2", "Notch3", "M2", 3, 1, 1.562510666666667, 4.687532, 0.06247701111592723, 0.06247701111592723, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 0.17374337358, 1.56369036222, 0.00007263421010672428, 0.00007263421010672429)
    ,@("M2", "Jag2", "Notch3", "sCs", 3, 1, 1.562510666666667, 4.687532, 0.06247701111592723, 0.06247701111592723, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 111.9052102885205, 1007.146892596684, 0.04678248377853075, 0.04678248377853075)
    ,@("sCs", "Jag2", "Notch3", "ECs", 3, 1, 1.953751, 5.861253, 0.07812076138024482, 0.07812076138024483, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 33.84858628742333, 304.63727658681, 0.01415055594672362, 0.01415055594672362)
    ,@("sCs", "Jag2", "Notch3", "FAPs", 3, 1, 1.953751, 5.861253, 0.07812076138024482, 0.07812076138024483, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 12.87615623922967, 115.885406153067, 0.005382935869013626, 0.005382935869013627)
    ,@("sCs", "Jag2", "Notch3", "M2", 3, 1, 1.953751, 5.861253, 0.07812076138024482, 0.07812076138024483, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 0.217247342445, 1.955226082005, 0.00009082124279699168, 0.00009082124279699169)
    ,@("sCs", "Jag2", "Notch3", "sCs", 3, 1, 1.953751, 5.861253, 0.07812076138024482, 0.07812076138024483, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 139.9253913400957, 1259.328522060861, 0.05849644832171059, 0.0584964483217106)
)

$startRow = 2
$startCol = 1
for ($i = 0; $i -lt $data.Count; $i++) {
    $dataRow = $data[$i]
    for ($j = 0; $j -lt $dataRow.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $dataRow[$j]
    }
}
